# Update the validation rule text for the "Content" row on the Comments sheet
# (shared-string content change: required|min:3|max:256 -> required|minlength:3|maxlength:256)
$wb = $excel.ActiveWorkbook
$wsPosts = $wb.Worksheets.Item("Posts")
$wsComments = $wb.Worksheets.Item("Comments")

$wsComments.Range("C3").Value = "required|minlength:3|maxlength:256"

# Move the selection/active view from Posts (D18) to Posts(D12) then over to
# Comments (C4), which becomes the active tab.
$wsPosts.Range("D12").Select()
$wsComments.Activate()
$wsComments.Range("C4").Select()
